$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("I2").Value = 0.7376376588883126
$ws.Range("J2").Value = 0.7376376588883125
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 0.8921946666666667
$ws.Range("N2").Value = 2.676584
$ws.Range("O2").Value = 0.1932069095523364
$ws.Range("P2").Value = 0.1932069095523364
$ws.Range("Q2").Value = 0.5197700105351112
$ws.Range("R2").Value = 4.677930094816
$ws.Range("S2").Value = 0.1425166924432313
$ws.Range("T2").Value = 0.1425166924432313

# Row 3
$ws.Range("I3").Value = 0.7376376588883126
$ws.Range("J3").Value = 0.7376376588883125
$ws.Range("O3").Value = 0.427802038465628
$ws.Range("P3").Value = 0.427802038465628
$ws.Range("S3").Value = 0.3155628941214337
$ws.Range("T3").Value = 0.3155628941214336

# Row 4
$ws.Range("I4").Value = 0.7376376588883126
$ws.Range("J4").Value = 0.7376376588883125
$ws.Range("M4").Value = 1.750112333333333
$ws.Range("N4").Value = 5.250337
$ws.Range("O4").Value = 0.3789910519820356
$ws.Range("P4").Value = 0.3789910519820357
$ws.Range("Q4").Value = 1.019571109220889
$ws.Range("R4").Value = 9.176139982988
$ws.Range("S4").Value = 0.2795580723236475
$ws.Range("T4").Value = 0.2795580723236475

# Row 5
$ws.Range("E5").Value = 2
$ws.Range("F5").Value = 0.6666666666666666
$ws.Range("G5").Value = 0.2072096666666667
$ws.Range("H5").Value = 0.621629
$ws.Range("I5").Value = 0.2623623411116874
$ws.Range("J5").Value = 0.2623623411116874
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 0.8921946666666667
$ws.Range("N5").Value = 2.676584
$ws.Range("O5").Value = 0.1932069095523364
$ws.Range("P5").Value = 0.1932069095523364
$ws.Range("Q5").Value = 0.1848713594817778
$ws.Range("R5").Value = 1.663842235336
$ws.Range("S5").Value = 0.05069021710910501
$ws.Range("T5").Value = 0.05069021710910502

# Row 6
$ws.Range("E6").Value = 2
$ws.Range("F6").Value = 0.6666666666666666
$ws.Range("G6").Value = 0.2072096666666667
$ws.Range("H6").Value = 0.621629
$ws.Range("I6").Value = 0.2623623411116874
$ws.Range("J6").Value = 0.2623623411116874
$ws.Range("O6").Value = 0.427802038465628
$ws.Range("P6").Value = 0.427802038465628
$ws.Range("Q6").Value = 0.4093453211557778
$ws.Range("R6").Value = 3.684107890402
$ws.Range("S6").Value = 0.1122391443441943
$ws.Range("T6").Value = 0.1122391443441943

# Row 7
$ws.Range("E7").Value = 2
$ws.Range("F7").Value = 0.6666666666666666
$ws.Range("G7").Value = 0.2072096666666667
$ws.Range("H7").Value = 0.621629
$ws.Range("I7").Value = 0.2623623411116874
$ws.Range("J7").Value = 0.2623623411116874
$ws.Range("M7").Value = 1.750112333333333
$ws.Range("N7").Value = 5.250337
$ws.Range("O7").Value = 0.3789910519820356
$ws.Range("P7").Value = 0.3789910519820357
$ws.Range("Q7").Value = 0.3626401932192222
$ws.Range("R7").Value = 3.263761738973
$ws.Range("S7").Value = 0.0994329796583881
$ws.Range("T7").Value = 0.09943297965838811

